$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.065.60"
$ws.Range("E2").Value = "  -2.99%  "
$ws.Range("D3").Value = "3.451.20"
$ws.Range("E3").Value = "  -4.13%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'583.48"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").Value = "'165.22"
$ws.Range("E6").Value = "  -4.99%  "
$ws.Range("E7").Value = "  -4.03%  "
$ws.Range("D8").Value = "3.444.23"
$ws.Range("E8").Value = "  -4.14%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  -3.81%  "
$ws.Range("E11").Value = "  +2.50%  "
$ws.Range("D12").Value = "'0.567"
$ws.Range("E12").Value = "  -7.56%  "
$ws.Range("D13").Value = "'45.83"
$ws.Range("E13").Value = "  -4.58%  "
$ws.Range("D14").Value = "'0.0000269"
$ws.Range("E14").Value = "  -3.64%  "
$ws.Range("D15").Value = "4.003.84"
$ws.Range("E15").Value = "  -4.25%  "
$ws.Range("D16").Value = "'609.81"
$ws.Range("E16").Value = "  -10.78%  "
$ws.Range("E17").Value = "  -7.91%  "
$ws.Range("D18").Value = "3.453.79"
$ws.Range("E18").Value = "  -4.10%  "
$ws.Range("D19").Value = "68.189.75"
$ws.Range("E19").Value = "  -2.92%  "
$ws.Range("E20").Value = "  -3.34%  "
$ws.Range("D21").Value = "'17.08"
$ws.Range("E21").Value = "  -3.12%  "
$ws.Range("D22").Value = "'10.90"
$ws.Range("E22").Value = "  -3.96%  "
$ws.Range("D23").Value = "'0.865"
$ws.Range("E23").Value = "  -6.90%  "
$ws.Range("D24").Value = "'15.49"
$ws.Range("E24").Value = "  -8.73%  "
$ws.Range("D25").Value = "'95.03"
$ws.Range("E25").Value = "  -3.95%  "
$ws.Range("D26").Value = "'3.72"
$ws.Range("E26").Value = "  -4.10%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -6.72%  "
$ws.Range("D29").Value = "'8.99"
$ws.Range("E29").Value = "  -6.73%  "
$ws.Range("D30").Value = "'32.19"
$ws.Range("E30").Value = "  -5.78%  "
$ws.Range("D31").Value = "'8.30"
$ws.Range("E31").Value = "  -8.67%  "
$ws.Range("D32").Value = "'3.03"
$ws.Range("E32").Value = "  -5.99%  "
$ws.Range("E33").Value = "  -5.91%  "
$ws.Range("D34").Value = "'6.73"
$ws.Range("E34").Value = "  -10.02%  "
$ws.Range("D35").Value = "'579.60"
$ws.Range("E35").Value = "  +1.68%  "
$ws.Range("E36").Value = "  -3.69%  "
$ws.Range("D37").Value = "'56.75"
$ws.Range("E37").Value = "  -2.64%  "
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("D39").Value = "'0.100"
$ws.Range("E39").Value = "  -5.97%  "
$ws.Range("D40").Value = "'3.37"
$ws.Range("E40").Value = "  -14.64%  "
$ws.Range("E41").Value = "  -3.22%  "
$ws.Range("D42").Value = "'0.0430"
$ws.Range("D43").Value = "3.357.98"
$ws.Range("E43").Value = "  -4.33%  "
$ws.Range("E44").Value = "  -7.54%  "
$ws.Range("D45").Value = "'32.34"
$ws.Range("E45").Value = "  -4.85%  "
$ws.Range("E46").Value = "  -5.88%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "'2.73"
$ws.Range("E47").Value = "  -7.89%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "'2.47"
$ws.Range("E48").Value = "  -5.84%  "
$ws.Range("E49").Value = "  -4.97%  "
$ws.Range("D50").Value = "'131.78"
$ws.Range("E50").Value = "  -2.79%  "
$ws.Range("D51").Value = "'5.57"
$ws.Range("E51").Value = "  +9.97%  "
